$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (Employees -> Sheet1) ---
$ws.Name = "Sheet1"

# --- Replace the employee sample data (rows 2-4) ---
# Columns: A=Name B=Emp.N C=ID.N D=Mobil E=Job Title F=Status G=Location H=Project I=Email
$data = @(
  @("موظف الأول",  "N-5001001", "EMP1001", "0599123456", "مدير قسم",      "نشط",      "الرياض", "مشروع التطوير", "employee1@example.com"),
  @("موظف الثاني", "N-5002002", "EMP1002", "0599234567", "مهندس برمجيات", "نشط",      "جدة",    "مشروع الدعم",   "employee2@example.com"),
  @("موظف الثالث", "N-5003003", "EMP1003", "0599345678", "محاسب",         "في إجازة", "الدمام", "مشروع المالية", "employee3@example.com")
)

# --- Strip the per-row highlight style from the data rows (A2:I4) first ---
# Previously centered/filled (style index 2); now plain default formatting.
$dataRange = $ws.Range("A2:I4")
$dataRange.Style = "Normal"

# Mobile numbers start with "0" - format the column as text so the
# leading zero survives (otherwise Excel coerces the digit string to a
# number and the zero is lost).
$ws.Range("D2:D4").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = 2 + $i
  $rowData = $data[$i]
  for ($c = 0; $c -lt $rowData.Length; $c++) {
    $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
  }
}

# --- Re-style the header row (A1:I1) ---
# Start from a clean slate so the old blue fill / white font don't linger,
# then rebuild: bold, automatic (default) text color, no fill, a thin box
# border around each header cell, centered horizontally and top-aligned.
$headerRange = $ws.Range("A1:I1")
$headerRange.Style = "Normal"
$headerRange.Font.Bold = $true
$headerRange.Font.ColorIndex = -4105
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
